$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph 2: +++IMAGE ({ width: 3, height: 3, path: './<...>/sample.png' })+++ ---
$xmlPng = @"
<w:p $wNs>
<w:r><w:t xml:space='preserve'>+++IMAGE </w:t></w:r>
<w:proofErr w:type='gramStart'/>
<w:r><w:t>(</w:t></w:r>
<w:r><w:t>{ width</w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:t>: 3, height: 3, path: '</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:r><w:t>/src/__tests__/fixtures/</w:t></w:r>
<w:r><w:t>sample</w:t></w:r>
<w:r><w:t>.png</w:t></w:r>
<w:r><w:t>' }</w:t></w:r>
<w:r><w:t>)</w:t></w:r>
<w:r><w:t>+++</w:t></w:r>
</w:p>
"@

# --- Paragraph 3: +++IMAGE ({ width: 3, height: 3, path: './<...>/sample.jpg' })+++ ---
$xmlJpg = @"
<w:p $wNs>
<w:r><w:t xml:space='preserve'>+++IMAGE </w:t></w:r>
<w:proofErr w:type='gramStart'/>
<w:r><w:t>({ width</w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:t>: 3, height: 3, path: '.</w:t></w:r>
<w:r><w:t>/src/__tests__/fixtures/</w:t></w:r>
<w:r><w:t>sample.jpg</w:t></w:r>
<w:r><w:t>'</w:t></w:r>
<w:r><w:t xml:space='preserve'> </w:t></w:r>
<w:r><w:t>})+++</w:t></w:r>
</w:p>
"@

# --- Paragraph 4: +++IMAGE ({ width: 3, height: 3, path: './<...>/sample.jpeg' })+++ ---
$xmlJpeg = @"
<w:p $wNs>
<w:r><w:t xml:space='preserve'>+++IMAGE </w:t></w:r>
<w:proofErr w:type='gramStart'/>
<w:r><w:t>({ width</w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:t>: 3, height: 3, path: './</w:t></w:r>
<w:r><w:t>src/__tests__/fixtures/</w:t></w:r>
<w:r><w:t>sample.jpeg' })+++</w:t></w:r>
</w:p>
"@

# --- Paragraph 5: +++IMAGE ({ width: 3, height: 3, path: './<...>/sample.gif' })+++ (keeps _GoBack bookmark) ---
$xmlGif = @"
<w:p $wNs>
<w:r><w:t xml:space='preserve'>+++IMAGE </w:t></w:r>
<w:proofErr w:type='gramStart'/>
<w:r><w:t>({ width</w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:t>: 3, height: 3, path: './</w:t></w:r>
<w:bookmarkStart w:id='0' w:name='_GoBack'/>
<w:bookmarkEnd w:id='0'/>
<w:r><w:t>src/__tests__/fixtures/</w:t></w:r>
<w:r><w:t>sample.gif' })+++</w:t></w:r>
</w:p>
"@

# --- Paragraph 6: previously empty, now holds a left tab stop + a tab run ---
$xmlTab = @"
<w:p $wNs>
<w:pPr>
<w:tabs>
<w:tab w:val='left' w:pos='1478'/>
</w:tabs>
</w:pPr>
<w:r><w:tab/></w:r>
</w:p>
"@

$d.Paragraphs(2).Range.InsertXML($xmlPng)
$d.Paragraphs(3).Range.InsertXML($xmlJpg)
$d.Paragraphs(4).Range.InsertXML($xmlJpeg)
$d.Paragraphs(5).Range.InsertXML($xmlGif)
$d.Paragraphs(6).Range.InsertXML($xmlTab)

Write-Output "paragraphs updated"
